$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A122").Value = 121
$ws.Range("B122").Value = 1
$ws.Range("C122").Value = "2024-06-17 12:24:11"
$ws.Range("D122").Value = 200
$ws.Range("E122").Value = 13

$ws.Range("A123").Value = 122
$ws.Range("B123").Value = 2
$ws.Range("C123").Value = "2024-06-17 12:24:12"
$ws.Range("D123").Value = 200
$ws.Range("E123").Value = 0
